# Auto-generated script applying scheduled-runner market data refresh
# to the Leve profit tables (H:N columns) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 2899.5557
$ws.Cells.Item(29, 9).Value = 1750.75
$ws.Cells.Item(29, 10).Value = 3818.6
$ws.Cells.Item(29, 11).Value = 5252.25
$ws.Cells.Item(29, 12).Value = 11455.8
$ws.Cells.Item(29, 13).Value = -4971.25
$ws.Cells.Item(29, 14).Value = -12017.8
$ws.Cells.Item(41, 8).Value = 3198.7144
$ws.Cells.Item(41, 10).Value = 3224.125
$ws.Cells.Item(41, 12).Value = 3224.125
$ws.Cells.Item(41, 14).Value = -4104.125
$ws.Cells.Item(49, 8).Value = 5165.6665
$ws.Cells.Item(49, 9).Value = 0
$ws.Cells.Item(49, 10).Value = 5165.6665
$ws.Cells.Item(49, 11).Value = 0
$ws.Cells.Item(49, 12).Value = 15496.9995
$ws.Cells.Item(49, 13).ClearContents()
$ws.Cells.Item(49, 14).Value = -15768.9995
$ws.Cells.Item(76, 8).Value = 5374.375
$ws.Cells.Item(76, 9).Value = 4999
$ws.Cells.Item(76, 11).Value = 4999
$ws.Cells.Item(76, 13).Value = -4684
$ws.Cells.Item(79, 8).Value = 5374.375
$ws.Cells.Item(79, 9).Value = 4999
$ws.Cells.Item(79, 11).Value = 4999
$ws.Cells.Item(79, 13).Value = -3907
$ws.Cells.Item(86, 8).Value = 1704142.9
$ws.Cells.Item(86, 9).Value = 2311281.8
$ws.Cells.Item(86, 10).Value = 4153.6
$ws.Cells.Item(86, 11).Value = 2311281.8
$ws.Cells.Item(86, 12).Value = 4153.6
$ws.Cells.Item(86, 13).Value = -2310158.8
$ws.Cells.Item(86, 14).Value = -6399.6
$ws.Cells.Item(89, 8).Value = 1704142.9
$ws.Cells.Item(89, 9).Value = 2311281.8
$ws.Cells.Item(89, 10).Value = 4153.6
$ws.Cells.Item(89, 11).Value = 11556409
$ws.Cells.Item(89, 12).Value = 20768
$ws.Cells.Item(89, 13).Value = -11550793
$ws.Cells.Item(89, 14).Value = -32000
$ws.Cells.Item(132, 8).Value = 9351.831
$ws.Cells.Item(132, 9).Value = 2793.1072
$ws.Cells.Item(132, 11).Value = 8379.321599999999
$ws.Cells.Item(132, 13).Value = -5849.321599999999
$ws.Cells.Item(137, 8).Value = 9013034
$ws.Cells.Item(137, 10).Value = 13894686
$ws.Cells.Item(137, 12).Value = 41684058
$ws.Cells.Item(137, 14).Value = -41689158

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5212.0386
$ws.Cells.Item(32, 9).Value = 2664.9302
$ws.Cells.Item(32, 11).Value = 2664.9302
$ws.Cells.Item(32, 13).Value = -2377.9302
$ws.Cells.Item(45, 8).Value = 2924.261
$ws.Cells.Item(45, 9).Value = 2671.6843
$ws.Cells.Item(45, 11).Value = 2671.6843
$ws.Cells.Item(45, 13).Value = -2294.6843
$ws.Cells.Item(132, 8).Value = 4030.2766
$ws.Cells.Item(132, 9).Value = 1606.3715
$ws.Cells.Item(132, 11).Value = 4819.1145
$ws.Cells.Item(132, 13).Value = -2289.1145

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(87, 8).Value = 50000
$ws.Cells.Item(87, 9).Value = 50000
$ws.Cells.Item(87, 11).Value = 50000
$ws.Cells.Item(87, 13).Value = -48752
$ws.Cells.Item(90, 8).Value = 50000
$ws.Cells.Item(90, 9).Value = 50000
$ws.Cells.Item(90, 11).Value = 150000
$ws.Cells.Item(90, 13).Value = -143760
$ws.Cells.Item(107, 8).Value = 1461.619
$ws.Cells.Item(107, 9).Value = 1449.05
$ws.Cells.Item(107, 10).Value = 1713
$ws.Cells.Item(107, 11).Value = 1449.05
$ws.Cells.Item(107, 12).Value = 1713
$ws.Cells.Item(107, 13).Value = 470.95
$ws.Cells.Item(107, 14).Value = -5553
$ws.Cells.Item(134, 8).Value = 4600.933
$ws.Cells.Item(134, 9).Value = 1902.5
$ws.Cells.Item(134, 11).Value = 5707.5
$ws.Cells.Item(134, 13).Value = -3172.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 270.16666
$ws.Cells.Item(7, 9).Value = 166.66667
$ws.Cells.Item(7, 11).Value = 166.66667
$ws.Cells.Item(7, 13).Value = -53.66667000000001
$ws.Cells.Item(31, 8).Value = 4893.641
$ws.Cells.Item(31, 9).Value = 1721.1364
$ws.Cells.Item(31, 10).Value = 8999.235000000001
$ws.Cells.Item(31, 11).Value = 1721.1364
$ws.Cells.Item(31, 12).Value = 8999.235000000001
$ws.Cells.Item(31, 13).Value = -1426.1364
$ws.Cells.Item(31, 14).Value = -9589.235000000001
$ws.Cells.Item(34, 8).Value = 4893.641
$ws.Cells.Item(34, 9).Value = 1721.1364
$ws.Cells.Item(34, 10).Value = 8999.235000000001
$ws.Cells.Item(34, 11).Value = 1721.1364
$ws.Cells.Item(34, 12).Value = 8999.235000000001
$ws.Cells.Item(34, 13).Value = -1519.1364
$ws.Cells.Item(34, 14).Value = -9403.235000000001
$ws.Cells.Item(122, 8).Value = 3286.6
$ws.Cells.Item(122, 9).Value = 2082.25
$ws.Cells.Item(122, 10).Value = 5427.6665
$ws.Cells.Item(122, 11).Value = 6246.75
$ws.Cells.Item(122, 12).Value = 16282.9995
$ws.Cells.Item(122, 13).Value = -3796.75
$ws.Cells.Item(122, 14).Value = -21182.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 95.545456
$ws.Cells.Item(2, 9).Value = 56.875
$ws.Cells.Item(2, 11).Value = 341.25
$ws.Cells.Item(2, 13).Value = -228.25
$ws.Cells.Item(5, 8).Value = 870.7778
$ws.Cells.Item(5, 9).Value = 646.7143
$ws.Cells.Item(5, 11).Value = 1940.1429
$ws.Cells.Item(5, 13).Value = -1828.1429
$ws.Cells.Item(14, 8).Value = 279.0909
$ws.Cells.Item(14, 9).Value = 279.0909
$ws.Cells.Item(14, 11).Value = 837.2727
$ws.Cells.Item(14, 13).Value = -664.2727
$ws.Cells.Item(86, 8).Value = 330.8
$ws.Cells.Item(86, 9).Value = 364.33334
$ws.Cells.Item(86, 11).Value = 1093.00002
$ws.Cells.Item(86, 13).Value = 92.99998000000005
$ws.Cells.Item(87, 8).Value = 17500
$ws.Cells.Item(89, 8).Value = 330.8
$ws.Cells.Item(89, 9).Value = 364.33334
$ws.Cells.Item(89, 11).Value = 3279.00006
$ws.Cells.Item(89, 13).Value = 2648.99994
$ws.Cells.Item(90, 8).Value = 17500
$ws.Cells.Item(96, 8).Value = 12998.6
$ws.Cells.Item(96, 9).Value = 12996
$ws.Cells.Item(96, 10).Value = 12999.25
$ws.Cells.Item(96, 11).Value = 38988
$ws.Cells.Item(96, 12).Value = 38997.75
$ws.Cells.Item(96, 13).Value = -36929
$ws.Cells.Item(96, 14).Value = -43115.75
$ws.Cells.Item(97, 8).Value = 216.83333
$ws.Cells.Item(122, 8).Value = 782.875
$ws.Cells.Item(122, 9).Value = 724
$ws.Cells.Item(122, 10).Value = 818.2
$ws.Cells.Item(122, 11).Value = 6516
$ws.Cells.Item(122, 12).Value = 7363.8
$ws.Cells.Item(122, 13).Value = -4066
$ws.Cells.Item(122, 14).Value = -12263.8
$ws.Cells.Item(129, 8).Value = 2015
$ws.Cells.Item(129, 10).Value = 1691.6666
$ws.Cells.Item(129, 12).Value = 5074.9998
$ws.Cells.Item(129, 14).Value = -15074.9998
$ws.Cells.Item(131, 8).Value = 10002852
$ws.Cells.Item(131, 10).Value = 8133024
$ws.Cells.Item(131, 12).Value = 24399072
$ws.Cells.Item(131, 14).Value = -24409152
$ws.Cells.Item(135, 8).Value = 870.7778
$ws.Cells.Item(135, 9).Value = 646.7143
$ws.Cells.Item(135, 11).Value = 5820.428699999999
$ws.Cells.Item(135, 13).Value = -3285.428699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(93, 8).Value = 44147.855
$ws.Cells.Item(93, 10).Value = 44147.855
$ws.Cells.Item(93, 12).Value = 44147.855
$ws.Cells.Item(93, 14).Value = -47891.855
$ws.Cells.Item(97, 8).Value = 369.66666
$ws.Cells.Item(97, 9).Value = 362.5
$ws.Cells.Item(97, 11).Value = 362.5
$ws.Cells.Item(97, 13).Value = 133.5
$ws.Cells.Item(102, 8).Value = 14577505
$ws.Cells.Item(102, 9).Value = 18894548
$ws.Cells.Item(102, 11).Value = 18894548
$ws.Cells.Item(102, 13).Value = -18892926
$ws.Cells.Item(126, 8).Value = 3997
$ws.Cells.Item(126, 9).Value = 3073.3684
$ws.Cells.Item(126, 10).Value = 5166.933
$ws.Cells.Item(126, 11).Value = 9220.1052
$ws.Cells.Item(126, 12).Value = 15500.799
$ws.Cells.Item(126, 13).Value = -6750.1052
$ws.Cells.Item(126, 14).Value = -20440.799

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 6978.3213
$ws.Cells.Item(46, 9).Value = 4652.6665
$ws.Cells.Item(46, 10).Value = 7612.591
$ws.Cells.Item(46, 11).Value = 4652.6665
$ws.Cells.Item(46, 12).Value = 7612.591
$ws.Cells.Item(46, 13).Value = -4464.6665
$ws.Cells.Item(46, 14).Value = -7988.591
$ws.Cells.Item(93, 8).Value = 3139.1428
$ws.Cells.Item(93, 9).Value = 624.5
$ws.Cells.Item(93, 10).Value = 6492
$ws.Cells.Item(93, 11).Value = 624.5
$ws.Cells.Item(93, 12).Value = 6492
$ws.Cells.Item(93, 13).Value = 623.5
$ws.Cells.Item(93, 14).Value = -8988

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 24390
$ws.Cells.Item(41, 10).Value = 24390
$ws.Cells.Item(41, 12).Value = 24390
$ws.Cells.Item(41, 14).Value = -25170
$ws.Cells.Item(45, 8).Value = 32916.332
$ws.Cells.Item(45, 10).Value = 32916.332
$ws.Cells.Item(45, 12).Value = 32916.332
$ws.Cells.Item(45, 14).Value = -33898.332
$ws.Cells.Item(106, 8).Value = 99988.5
$ws.Cells.Item(106, 10).Value = 99988.5
$ws.Cells.Item(106, 12).Value = 99988.5
$ws.Cells.Item(106, 14).Value = -102512.5
$ws.Cells.Item(132, 8).Value = 26457008
$ws.Cells.Item(132, 9).Value = 3473594.5
$ws.Cells.Item(132, 10).Value = 100003940
$ws.Cells.Item(132, 11).Value = 10420783.5
$ws.Cells.Item(132, 12).Value = 300011820
$ws.Cells.Item(132, 13).Value = -10418253.5
$ws.Cells.Item(132, 14).Value = -300016880
$ws.Cells.Item(136, 8).Value = 7757.244
$ws.Cells.Item(136, 9).Value = 2890.3784
$ws.Cells.Item(136, 11).Value = 8671.135200000001
$ws.Cells.Item(136, 13).Value = -6121.135200000001

"Applied 209 cell updates."
